# Update workbook/sheet data to reflect revised figures and bookkeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "სოციალური პაკეტის მიმღებები" row (row 4) values for columns E:K.
$ws.Range("E4").Value = 6222
$ws.Range("F4").Value = 6138
$ws.Range("G4").Value = 5987
$ws.Range("H4").Value = 5864
$ws.Range("I4").Value = 6831
$ws.Range("J4").Value = 6994
$ws.Range("K4").Value = 7212

# Move the active selection to E4:K4 (active cell E4), matching the saved view state.
$ws.Range("E4:K4").Select()
